# Adds a new "localdb" command category (new-column) to the '#system' sheet,
# inserts it alphabetically into the "target" category list (column A), wires
# up the new named range, and shifts the named ranges that come after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. Make room: insert a brand-new column before N (macro's old column).
#        This shifts N:AC -> O:AD, matching the author's diff exactly.
$ws.Range("N1").EntireColumn.Insert()

# --- 2. Make room in the "target" list (column A only): shift A14:A29 down
#        to A15:A30 (a plain Range("A14").Insert() would shift the WHOLE row,
#        not just column A, so the values are moved by hand instead).
for ($r = 29; $r -ge 14; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
}

# --- 3. Populate the new column N with the "localdb" category header plus
#        its six function names (rows 2-7).
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# --- 4. Insert "localdb" into the alphabetical "target" list in column A.
$ws.Range("A14").Value = "localdb"

# --- 5. Update the defined (named) ranges that shifted one column to the
#        right because of the new column N.
$names = $wb.Names
$names.Item("macro").RefersTo     = "='#system'!`$O`$2:`$O`$4"
$names.Item("mail").RefersTo      = "='#system'!`$P`$2:`$P`$2"
$names.Item("number").RefersTo    = "='#system'!`$Q`$2:`$Q`$16"
$names.Item("pdf").RefersTo       = "='#system'!`$R`$2:`$R`$16"
$names.Item("rdbms").RefersTo     = "='#system'!`$S`$2:`$S`$7"
$names.Item("redis").RefersTo     = "='#system'!`$T`$2:`$T`$10"
$names.Item("sms").RefersTo       = "='#system'!`$U`$2:`$U`$2"
$names.Item("sound").RefersTo     = "='#system'!`$V`$2:`$V`$5"
$names.Item("ssh").RefersTo       = "='#system'!`$W`$2:`$W`$9"
$names.Item("step").RefersTo      = "='#system'!`$X`$2:`$X`$4"
$names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$127"
$names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$21"

# --- 6. The "target" list grew by one row (new "localdb" entry).
$names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"

# --- 7. Register the new "localdb" named range itself.
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")

# --- 8. The sheet's "dimension" hint has always extended one column past the
#        real data (was A1:AD127 while real data stopped at AC); touching the
#        new trailing column keeps that same one-column pad (now A1:AE127).
$ws.Range("AE1").Font.Bold = $false

"localdb command category added"
